# Update Name of Algo
# Applies updated imputed values to columns A and B of Sheet1,
# matching the values produced by the (re-)run of the RandomForest
# imputation algorithm.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6.700700000000007
$ws.Range("B4").Value = 9.003999999999998
$ws.Range("A11").Value = -21.68189999999999
$ws.Range("A12").Value = -21.55560000000001
$ws.Range("B14").Value = 6.782999999999995
$ws.Range("A15").Value = -21.90789999999999
$ws.Range("B26").Value = 4.410000000000002
$ws.Range("A27").Value = -21.70319999999999
$ws.Range("A28").Value = -21.83449999999999
$ws.Range("A31").Value = -21.60369999999999
$ws.Range("B31").Value = 4.671200000000003
$ws.Range("A32").Value = -21.363
$ws.Range("B35").Value = 9.113000000000007
$ws.Range("A36").Value = -19.73389999999999
$ws.Range("B37").Value = 8.925200000000004
$ws.Range("A38").Value = -19.2326
$ws.Range("B39").Value = 9.521600000000007
$ws.Range("B40").Value = 8.674200000000004
$ws.Range("B45").Value = 6.461499999999996
$ws.Range("A46").Value = -21.6771
$ws.Range("B52").Value = 5.186900000000001
$ws.Range("A54").Value = -21.66229999999999
$ws.Range("A55").Value = -22.3875
$ws.Range("A56").Value = -21.9934
$ws.Range("B57").Value = 4.628999999999996
$ws.Range("A67").Value = -21.52299999999997
$ws.Range("A69").Value = -21.72669999999997
$ws.Range("A72").Value = -22.18100000000003
$ws.Range("A73").Value = -19.87139999999998
$ws.Range("B81").Value = 6.2989
$ws.Range("A83").Value = -21.61859999999999
$ws.Range("B83").Value = 5.671100000000004
$ws.Range("A86").Value = -22.0123
$ws.Range("A91").Value = -21.4221
$ws.Range("A93").Value = -21.30209999999999
$ws.Range("A99").Value = -20.21759999999999
$ws.Range("B100").Value = 5.390799999999996
$ws.Range("B102").Value = 8.194200000000002
